# Apply the recorded edit:
#  - On the "role" sheet, clear the contents of B29:I35 (values/types), which
#    were sample "Camp_Light"/"Camp_Dark" demo rows. Clearing these also
#    removes the now-orphaned "Camp_Light" shared string and the cell
#    comments that were attached to E29:G35.
#  - Explicitly delete the leftover comments on E29:G35 so they don't linger
#    attached to now-empty cells.
#  - Make "role" the active/selected sheet (tabSelected) with the cursor on
#    F21, matching the new workbookView activeTab="1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("role")

# Remove the cell comments that sit on E29:G35 ("光明"/"黑暗" notes) before
# the underlying cells are cleared.
foreach ($col in @("E", "F", "G")) {
    foreach ($row in 29..35) {
        $cell = $ws.Range("$col$row")
        if ($cell.Comment -ne $null) {
            $cell.Comment.Delete()
        }
    }
}

# Clear the sample data in columns B:I for rows 29-35 (keeps cell styles).
$ws.Range("B29:I35").ClearContents()

# Switch the active tab to "role" and set the new selection.
$ws.Activate()
$ws.Range("F21").Select()
